$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_3_6_0"
$ws.Range("B2").Value = 0.9871737183827846
$ws.Range("C2").Value = 0.9991490911100659
$ws.Range("D2").Value = 0.9809597415040151
$ws.Range("E2").Value = 0.9904013999521859
$ws.Range("F2").Value = 1.546853381812864
$ws.Range("G2").Value = 0.1256792613170171
$ws.Range("H2").Value = 2.949969822688463
$ws.Range("I2").Value = 1.45475715817014

$ws.Range("A3").Value = "model_3_6_1"
$ws.Range("B3").Value = 0.9887885310499448
$ws.Range("C3").Value = 0.9988947983168163
$ws.Range("D3").Value = 0.9805216335289381
$ws.Range("E3").Value = 0.9900594455530949
$ws.Range("F3").Value = 1.352106493374172
$ws.Range("G3").Value = 0.1632383123410571
$ws.Range("H3").Value = 3.017847331065191
$ws.Range("I3").Value = 1.506583529450048

$ws.Range("A4").Value = "model_3_6_2"
$ws.Range("B4").Value = 0.9900282292827653
$ws.Range("C4").Value = 0.9985693198862
$ws.Range("D4").Value = 0.9799294331660653
$ws.Range("E4").Value = 0.9896066378210326
$ws.Range("F4").Value = 1.202598517399905
$ws.Range("G4").Value = 0.211311483532907
$ws.Range("H4").Value = 3.109598879492348
$ws.Range("I4").Value = 1.575210754900745

$ws.Range("A5").Value = "model_3_6_3"
$ws.Range("B5").Value = 0.9909701874698097
$ws.Range("C5").Value = 0.9981922734474181
$ws.Range("D5").Value = 0.9792371821506851
$ws.Range("E5").Value = 0.9890790867927607
$ws.Range("F5").Value = 1.08899807959257
$ws.Range("G5").Value = 0.2670012506382831
$ws.Range("H5").Value = 3.216851604318939
$ws.Range("I5").Value = 1.655166022424715

$ws.Range("A6").Value = "model_3_6_4"
$ws.Range("B6").Value = 0.9916766871822777
$ws.Range("C6").Value = 0.9977803366350629
$ws.Range("D6").Value = 0.9784860769157948
$ws.Range("E6").Value = 0.9885052296327819
$ws.Range("F6").Value = 1.003794003922323
$ws.Range("G6").Value = 0.3278443266697126
$ws.Range("H6").Value = 3.333222806792742
$ws.Range("I6").Value = 1.74213941511613

$ws.Range("A7").Value = "model_3_6_5"
$ws.Range("B7").Value = 0.9921976411124165
$ws.Range("C7").Value = 0.9973474871241144
$ws.Range("D7").Value = 0.977706846016021
$ws.Range("E7").Value = 0.9879070493003157
$ws.Range("F7").Value = 0.9409668048436556
$ws.Range("G7").Value = 0.3917762087325037
$ws.Range("H7").Value = 3.453951610959127
$ws.Range("I7").Value = 1.832799210940207

$ws.Range("A8").Value = "model_3_6_24"
$ws.Range("B8").Value = 0.9924854766174188
$ws.Range("C8").Value = 0.9916371950011101
$ws.Range("D8").Value = 0.9683402156032379
$ws.Range("E8").Value = 0.9804549745719783
$ws.Range("F8").Value = 0.9062537572429366
$ws.Range("G8").Value = 1.235186477931921
$ws.Range("H8").Value = 4.905154443305794
$ws.Range("I8").Value = 2.962230482194864

$ws.Range("A9").Value = "model_3_6_23"
$ws.Range("B9").Value = 0.9925381023381767
$ws.Range("C9").Value = 0.9917761454701648
$ws.Range("D9").Value = 0.968552188284035
$ws.Range("E9").Value = 0.9806286436027983
$ws.Range("F9").Value = 0.8999070796512395
$ws.Range("G9").Value = 1.214663490668512
$ws.Range("H9").Value = 4.872312819242889
$ws.Range("I9").Value = 2.935909324476096

$ws.Range("A10").Value = "model_3_6_6"
$ws.Range("B10").Value = 0.9925731631974184
$ws.Range("C10").Value = 0.9969052252555747
$ws.Range("D10").Value = 0.9769222882941795
$ws.Range("E10").Value = 0.9873014435861298
$ws.Range("F10").Value = 0.8956787295880035
$ws.Range("G10").Value = 0.4570982962136246
$ws.Range("H10").Value = 3.575505717174517
$ws.Range("I10").Value = 1.924584392461685

$ws.Range("A11").Value = "model_3_6_22"
$ws.Range("B11").Value = 0.9925941702511226
$ws.Range("C11").Value = 0.9919287866909369
$ws.Range("D11").Value = 0.9687856695537059
$ws.Range("E11").Value = 0.9808197011916862
$ws.Range("F11").Value = 0.8931452726568324
$ws.Range("G11").Value = 1.192118379082412
$ws.Range("H11").Value = 4.836138798819916
$ws.Range("I11").Value = 2.906952768970825

$ws.Range("A12").Value = "model_3_6_21"
$ws.Range("B12").Value = 0.9926535068087485
$ws.Range("C12").Value = 0.9920961183193838
$ws.Range("D12").Value = 0.9690426585179317
$ws.Range("E12").Value = 0.9810296894908711
$ws.Range("F12").Value = 0.8859892661408487
$ws.Range("G12").Value = 1.167403493967257
$ws.Range("H12").Value = 4.796322654023882
$ws.Range("I12").Value = 2.875127088158052

$ws.Range("A13").Value = "model_3_6_20"
$ws.Range("B13").Value = 0.9927157324130648
$ws.Range("C13").Value = 0.9922793411164543
$ws.Range("D13").Value = 0.9693249853149395
$ws.Range("E13").Value = 0.9812600066975533
$ws.Range("F13").Value = 0.8784848397338344
$ws.Range("G13").Value = 1.140341482905132
$ws.Range("H13").Value = 4.752580835524689
$ws.Range("I13").Value = 2.840220372240974

$ws.Range("A14").Value = "model_3_6_19"
$ws.Range("B14").Value = 0.9927803556701876
$ws.Range("C14").Value = 0.9924795371259226
$ws.Range("D14").Value = 0.9696348955913168
$ws.Range("E14").Value = 0.9815123908268281
$ws.Range("F14").Value = 0.8706912556844986
$ws.Range("G14").Value = 1.110772528007351
$ws.Range("H14").Value = 4.704565417916413
$ws.Range("I14").Value = 2.801969208858595

$ws.Range("A15").Value = "model_3_6_7"
$ws.Range("B15").Value = 0.9928352302963068
$ws.Range("C15").Value = 0.9964626434689338
$ws.Range("D15").Value = 0.9761490556330125
$ws.Range("E15").Value = 0.9867011327084574
$ws.Range("F15").Value = 0.8640733594366631
$ws.Range("G15").Value = 0.5224676356051803
$ws.Range("H15").Value = 3.695305194520942
$ws.Range("I15").Value = 2.015567092237817

$ws.Range("A16").Value = "model_3_6_18"
$ws.Range("B16").Value = 0.9928465648682309
$ws.Range("C16").Value = 0.9926978506122565
$ws.Range("D16").Value = 0.9699745289490009
$ws.Range("E16").Value = 0.9817884130160979
$ws.Range("F16").Value = 0.8627064066879733
$ws.Range("G16").Value = 1.078527621387436
$ws.Range("H16").Value = 4.651944905639341
$ws.Range("I16").Value = 2.760135477517178

$ws.Range("A17").Value = "model_3_6_17"
$ws.Range("B17").Value = 0.9929132210037427
$ws.Range("C17").Value = 0.9929353751333538
$ws.Range("D17").Value = 0.9703459460325123
$ws.Range("E17").Value = 0.9820896453402316
$ws.Range("F17").Value = 0.8546676568996838
$ws.Range("G17").Value = 1.043445244520383
$ws.Range("H17").Value = 4.594400036265827
$ws.Range("I17").Value = 2.714480915641197

$ws.Range("A18").Value = "model_3_6_16"
$ws.Range("B18").Value = 0.9929788129196145
$ws.Range("C18").Value = 0.9931931498626714
$ws.Range("D18").Value = 0.9707515892984707
$ws.Range("E18").Value = 0.9824177686573233
$ws.Range("F18").Value = 0.8467572523168105
$ws.Range("G18").Value = 1.005371911464329
$ws.Range("H18").Value = 4.53155239196486
$ws.Range("I18").Value = 2.664750773545065

$ws.Range("A19").Value = "model_3_6_8"
$ws.Range("B19").Value = 0.9930094431671207
$ws.Range("C19").Value = 0.9960269323913573
$ws.Range("D19").Value = 0.975398868675515
$ws.Range("E19").Value = 0.9861154470549687
$ws.Range("F19").Value = 0.8430632353480159
$ws.Range("G19").Value = 0.5868221711203666
$ws.Range("H19").Value = 3.811534125260456
$ws.Range("I19").Value = 2.104333203192116

$ws.Range("A20").Value = "model_3_6_15"
$ws.Range("B20").Value = 0.9930411756044645
$ws.Range("C20").Value = 0.9934720617746284
$ws.Range("D20").Value = 0.9711931760130853
$ws.Range("E20").Value = 0.9827740973523257
$ws.Range("F20").Value = 0.8392362939566319
$ws.Range("G20").Value = 0.9641766160784856
$ws.Range("H20").Value = 4.463135911039046
$ws.Range("I20").Value = 2.610745843958043

$ws.Range("A21").Value = "model_3_6_14"
$ws.Range("B21").Value = 0.9930975737668394
$ws.Range("C21").Value = 0.9937727724494109
$ws.Range("D21").Value = 0.9716729554977765
$ws.Range("E21").Value = 0.9831600690693386
$ws.Range("F21").Value = 0.8324346587827591
$ws.Range("G21").Value = 0.9197616429551879
$ws.Range("H21").Value = 4.388802098728527
$ws.Range("I21").Value = 2.552248238538633

$ws.Range("A22").Value = "model_3_6_9"
$ws.Range("B22").Value = 0.9931161402316406
$ws.Range("C22").Value = 0.9956035105593105
$ws.Range("D22").Value = 0.9746798206475422
$ws.Range("E22").Value = 0.9855510738510884
$ws.Range("F22").Value = 0.8301955375999277
$ws.Range("G22").Value = 0.6493615848068932
$ws.Range("H22").Value = 3.922938599313715
$ws.Range("I22").Value = 2.189869214082692

$ws.Range("A23").Value = "model_3_6_13"
$ws.Range("B23").Value = 0.9931443059427271
$ws.Range("C23").Value = 0.9940957402418213
$ws.Range("D23").Value = 0.972192397499169
$ws.Range("E23").Value = 0.9835765633669226
$ws.Range("F23").Value = 0.8267987444571061
$ws.Range("G23").Value = 0.8720592930802538
$ws.Range("H23").Value = 4.308323242358603
$ws.Range("I23").Value = 2.489124652002149

$ws.Range("A24").Value = "model_3_6_10"
$ws.Range("B24").Value = 0.9931715092718052
$ws.Range("C24").Value = 0.9951963505655944
$ws.Range("D24").Value = 0.9739973169273848
$ws.Range("E24").Value = 0.985012683253062
$ws.Range("F24").Value = 0.8235180148710143
$ws.Range("G24").Value = 0.7094991246227398
$ws.Range("H24").Value = 4.028681143657923
$ws.Range("I24").Value = 2.271467319271875

$ws.Range("A25").Value = "model_3_6_12"
$ws.Range("B25").Value = 0.9931766030283987
$ws.Range("C25").Value = 0.994441045782996
$ws.Range("D25").Value = 0.972752672842302
$ws.Range("E25").Value = 0.9840242434460844
$ws.Range("F25").Value = 0.8229037063092718
$ws.Range("G25").Value = 0.8210576572331264
$ws.Range("H25").Value = 4.221517942158107
$ws.Range("I25").Value = 2.421274570064497

$ws.Range("A26").Value = "model_3_6_11"
$ws.Range("B26").Value = 0.9931882859834139
$ws.Range("C26").Value = 0.9948082407383708
$ws.Range("D26").Value = 0.9733543207167733
$ws.Range("E26").Value = 0.9845031255369446
$ws.Range("F26").Value = 0.821494738456064
$ws.Range("G26").Value = 0.76682295443855
$ws.Range("H26").Value = 4.1283026597107
$ws.Range("I26").Value = 2.3486955328999

